$runs = @(
    @{Text='editTask.html '; Bold=$true; Italic=$false},
    @{Text='-> hernoemd naar task.html (de naamgevings conventie is niet meer uniform, heb deze aangepast maar de rest zo gelaten) en css error gefixt en styling voorzien.'; Bold=$false; Italic=$false},
    @{Text=' Index Pagina'; Bold=$true; Italic=$false},
    @{Text=': naam owner weergegeven op klusje, bieden op klusje knop styling.  Zonder succes gezocht naar systeem om user names weer te geven op index pagina (join queries, meerdere queries combineren, etc). Count(offers) doet ook moeilijk, staat in comentaar in '; Bold=$false; Italic=$false},
    @{Text='taskrepo'; Bold=$false; Italic=$true},
    @{Text='. Toegewezen klusjes zijn ook van de homepagina afgehaald, hebben geen toegevoegde waarde van daar nog te staan. '; Bold=$false; Italic=$false},
    @{Text='Profile pagina'; Bold=$true; Italic=$false},
    @{Text=': Helemaal herstructureerd met nieuwe layout. Klusjesman en klant krijgen nu elk de juist taken op de juiste plaats te zien en kunnen hun vereiste acties uitvoeren op de taken ('; Bold=$false; Italic=$false},
    @{Text='/completeTask'; Bold=$true; Italic=$false},
    @{Text=' en '; Bold=$false; Italic=$false},
    @{Text='/finalizeTask'; Bold=$true; Italic=$false},
    @{Text=') '; Bold=$false; Italic=$false},
    @{Text='[maincontroller,taskcontroller] '; Bold=$false; Italic=$true},
    @{Text='en'; Bold=$false; Italic=$false},
    @{Text=' taskrepo'; Bold=$false; Italic=$true},
    @{Text='. Volledig proces is functioneel van taak aanmaken tot het beoordelen, inclusief css styling, '; Bold=$false; Italic=$false},
    @{Text='To do'; Bold=$true; Italic=$false},
    @{Text=': Rating berekenen en ophalen voor Profile page klusjesman als pagina toekennen van task bij klant. '; Bold=$false; Italic=$false},
    @{Text='Info Pagina:'; Bold=$true; Italic=$false},
    @{Text=' Informatie pagina gemaakt. Hier is de werking van de website uitgelegd voor zowel klant als klusjesman met visuele begeleiding. Project is bijna klaar!! (denk ik)'; Bold=$false; Italic=$false},
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full log entry text for row 15 from the formatted runs above
$fullText = ""
foreach ($r in $runs) { $fullText += $r.Text }

$cell = $ws.Range("A15")
$cell.Value = $fullText

# Apply per-run character-level formatting (bold / italic) so the shared
# string ends up with the same run structure as the source edit.
$pos = 1
foreach ($r in $runs) {
    $len = $r.Text.Length
    $chars = $cell.Characters($pos, $len)
    $chars.Font.Bold = $r.Bold
    $chars.Font.Italic = $r.Italic
    $pos += $len
}

# Date worked and hours for the new logbook entry
$ws.Range("B15").Value = 45642
$ws.Range("C15").Value = 10.5

# The long wrapped entry needs a much taller row; match the autofit height
# Excel computed for this text on the real worksheet.
$ws.Rows(15).RowHeight = 408.6

# Update the view to reflect where the user ended up after adding the row
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
